$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 26 with data (2020-04-19), matching the style of the row above (A25)
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A26").Value = 43940

$ws.Range("B26").Value = 48527
$ws.Range("C26").Value = 11555

$ws.Range("D26").Formula = "=B26-B25"
$ws.Range("E26").Formula = "=C26-C25"
$ws.Range("F26").Formula = "=E26/D26"

# Update view / selection to match diff
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("F27").Select()
